# The post at row 719 ("「あらゆる問題は、食べ物が解決してくれる」") was removed.
# Deleting the entire row shifts every subsequent row up by one and
# automatically shrinks the sheet's used range from A1:C769 to A1:C768.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("719").Delete()
